$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.7748443491049954
$ws.Range("D2").Value = 0.004241245969367924
$ws.Range("E2").Value = 0.4807100932212052
$ws.Range("F2").Value = 0.5738872358297726
$ws.Range("G2").Value = 0.002367532321094571
$ws.Range("I2").Value = 0.7740553189655508
$ws.Range("L2").Value = 0.2133349440246661
$ws.Range("M2").Value = 0.2090503303009541
$ws.Range("N2").Value = 1.574741273958978
$ws.Range("O2").Value = 1.839629889463225
# Row 3
$ws.Range("B3").Value = 0.7303581576005342
$ws.Range("D3").Value = 0.003887375460227105
$ws.Range("E3").Value = 0.4675654428869507
$ws.Range("F3").Value = 0.5506324538694258
$ws.Range("G3").Value = 0.002370847008024264
$ws.Range("I3").Value = 0.786366584358291
$ws.Range("L3").Value = 0.1956957426630623
$ws.Range("M3").Value = 0.1950314830238824
$ws.Range("N3").Value = 1.566145445541849
$ws.Range("O3").Value = 1.773912749601266
# Row 4
$ws.Range("B4").Value = 0.7032020066927203
$ws.Range("D4").Value = 0.003668921130433489
$ws.Range("E4").Value = 0.4594787765513644
$ws.Range("F4").Value = 0.5367620818859251
$ws.Range("G4").Value = 0.002372993154237461
$ws.Range("I4").Value = 0.7943611820636409
$ws.Range("L4").Value = 0.1849108055239412
$ws.Range("M4").Value = 0.1864688781993848
$ws.Range("N4").Value = 1.561434346395728
$ws.Range("O4").Value = 1.734905112674625
# Row 5
$ws.Range("B5").Value = 0.692176493205551
$ws.Range("D5").Value = 0.003579609102647652
$ws.Range("E5").Value = 0.456179752841031
$ws.Range("F5").Value = 0.5312121717730633
$ws.Range("G5").Value = 0.002373895702896095
$ws.Range("I5").Value = 0.797728529855922
$ws.Range("L5").Value = 0.1805275510443778
$ws.Range("M5").Value = 0.1829911269549243
$ws.Range("N5").Value = 1.559657839741647
$ws.Range("O5").Value = 1.71934596061385
# Row 6
$ws.Range("B6").Value = 0.6903482102040357
$ws.Range("D6").Value = 0.003564761554706308
$ws.Range("E6").Value = 0.4556317422753082
$ws.Range("F6").Value = 0.5302967894218682
$ws.Range("G6").Value = 0.00237404726281364
$ws.Range("I6").Value = 0.7982942877365193
$ws.Range("L6").Value = 0.1798004280365006
$ws.Range("M6").Value = 0.1824143560021625
$ws.Range("N6").Value = 1.559371530587867
$ws.Range("O6").Value = 1.716782688932057
# Row 7
$ws.Range("B7").Value = 0.7030531461237501
$ws.Range("D7").Value = 0.003667717804038517
$ws.Range("E7").Value = 0.4594342990664018
$ws.Range("F7").Value = 0.5366868196685886
$ws.Range("G7").Value = 0.002373005212952595
$ws.Range("I7").Value = 0.7944061520064023
$ws.Range("L7").Value = 0.1848516437670185
$ws.Range("M7").Value = 0.186421928801991
$ws.Range("N7").Value = 1.561409806526385
$ws.Range("O7").Value = 1.734693914137324
# Row 8
$ws.Range("B8").Value = 0.7594731695427583
$ws.Range("D8").Value = 0.004119478889016648
$ws.Range("E8").Value = 0.476181290963801
$ws.Range("F8").Value = 0.5657841114829836
$ws.Range("G8").Value = 0.002368652256964935
$ws.Range("I8").Value = 0.7782098658611034
$ws.Range("L8").Value = 0.207243613076912
$ws.Range("M8").Value = 0.20420743714228
$ws.Range("N8").Value = 1.571660260979655
$ws.Range("O8").Value = 1.816691311576676
# Row 9
$ws.Range("B9").Value = 0.8713338756544999
$ws.Range("D9").Value = 0.004995837596599273
$ws.Range("E9").Value = 0.5088837487412547
$ws.Range("F9").Value = 0.6260975740445787
$ws.Range("G9").Value = 0.002360992243782305
$ws.Range("I9").Value = 0.7499038805780929
$ws.Range("L9").Value = 0.2515083149291826
$ws.Range("M9").Value = 0.2394324833578665
$ws.Range("N9").Value = 1.596226618815635
$ws.Range("O9").Value = 1.988195915065887
# Row 10
$ws.Range("B10").Value = 0.9542214552155315
$ws.Range("D10").Value = 0.005633660349293024
$ws.Range("E10").Value = 0.5328107246589155
$ws.Range("F10").Value = 0.672419811898564
$ws.Range("G10").Value = 0.002355893032632402
$ws.Range("I10").Value = 0.731212934187941
$ws.Range("L10").Value = 0.2842384328994285
$ws.Range("M10").Value = 0.2655141661284475
$ws.Range("N10").Value = 1.616958675674226
$ws.Range("O10").Value = 2.120815983403304
# Row 11
$ws.Range("B11").Value = 0.9920737779063984
$ws.Range("D11").Value = 0.005922468826550897
$ws.Range("E11").Value = 0.5436710432402236
$ws.Range("F11").Value = 0.6939354085729548
$ws.Range("G11").Value = 0.002353686878641534
$ws.Range("I11").Value = 0.7231669498702749
$ws.Range("L11").Value = 0.2991722738750866
$ws.Range("M11").Value = 0.2774212433569616
$ws.Range("N11").Value = 1.62696493235174
$ws.Range("O11").Value = 2.182604589834568
# Row 12
$ws.Range("B12").Value = 1.006427643068776
$ws.Range("D12").Value = 0.00603163546118779
$ws.Range("E12").Value = 0.5477797790485681
$ws.Range("F12").Value = 0.7021469264392124
$ws.Range("G12").Value = 0.002352867698245577
$ws.Range("I12").Value = 0.7201858166923283
$ws.Range("L12").Value = 0.3048335812245284
$ws.Range("M12").Value = 0.2819360132832927
$ws.Range("N12").Value = 1.630836028985229
$ws.Range("O12").Value = 2.206213396308044
# Row 13
$ws.Range("B13").Value = 1.003335407751251
$ws.Range("D13").Value = 0.006008133415480899
$ws.Range("E13").Value = 0.546895065286364
$ws.Range("F13").Value = 0.7003755771590079
$ws.Range("G13").Value = 0.002353043402237565
$ws.Range("I13").Value = 0.7208249340298907
$ws.Range("L13").Value = 0.3036140462282049
$ws.Range("M13").Value = 0.28096342244352
$ws.Range("N13").Value = 1.629998686586717
$ws.Range("O13").Value = 2.201119426199455
# Row 14
$ws.Range("B14").Value = 0.9932542838655536
$ws.Range("D14").Value = 0.005931454052227281
$ws.Range("E14").Value = 0.5440091502106839
$ws.Range("F14").Value = 0.6946096903172077
$ws.Range("G14").Value = 0.002353619159023709
$ws.Range("I14").Value = 0.7229203731354432
$ws.Range("L14").Value = 0.2996379104861262
$ws.Range("M14").Value = 0.2777925611770442
$ws.Range("N14").Value = 1.627281771841012
$ws.Range("O14").Value = 2.184542670102701
# Row 15
$ws.Range("B15").Value = 0.9870818799627159
$ws.Range("D15").Value = 0.005884459648502371
$ws.Range("E15").Value = 0.5422409344969452
$ws.Range("F15").Value = 0.6910862645496678
$ws.Range("G15").Value = 0.002353973939626503
$ws.Range("I15").Value = 0.7242124483530588
$ws.Range("L15").Value = 0.2972032107483358
$ws.Range("M15").Value = 0.2758510665752567
$ws.Range("N15").Value = 1.62562823225835
$ws.Range("O15").Value = 2.174416410345032
# Row 16
$ws.Range("B16").Value = 0.951750573464011
$ws.Range("D16").Value = 0.00561475854087945
$ws.Range("E16").Value = 0.5321004604146111
$ws.Range("F16").Value = 0.6710226621590749
$ws.Range("G16").Value = 0.002356039487004865
$ws.Range("I16").Value = 0.73174795039799
$ws.Range("L16").Value = 0.2832633504211799
$ws.Range("M16").Value = 0.2647368398875258
$ws.Range("N16").Value = 1.616316251414645
$ws.Range("O16").Value = 2.11680739755036
# Row 17
$ws.Range("B17").Value = 0.9301127088241969
$ws.Range("D17").Value = 0.005448957987418623
$ws.Range("E17").Value = 0.525873162142581
$ws.Range("F17").Value = 0.6588280450048956
$ws.Range("G17").Value = 0.002357335647632181
$ws.Range("I17").Value = 0.7364877202494244
$ws.Range("L17").Value = 0.2747229812142962
$ws.Range("M17").Value = 0.2579292811685718
$ws.Range("N17").Value = 1.610750416786047
$ws.Range("O17").Value = 2.081840542257169
# Row 18
$ws.Range("B18").Value = 0.9176810177959283
$ws.Range("D18").Value = 0.005353468103695747
$ws.Range("E18").Value = 0.5222891297532968
$ws.Range("F18").Value = 0.6518557404777994
$ws.Range("G18").Value = 0.002358091853052401
$ws.Range("I18").Value = 0.7392568871188447
$ws.Range("L18").Value = 0.2698150132854096
$ws.Range("M18").Value = 0.254017759832017
$ws.Range("N18").Value = 1.607603275775332
$ws.Range("O18").Value = 2.061865807688491
# Row 19
$ws.Range("B19").Value = 0.9134742687432436
$ws.Range("D19").Value = 0.005321115452574077
$ws.Range("E19").Value = 0.52107526144421
$ws.Range("F19").Value = 0.649502194578119
$ws.Range("G19").Value = 0.002358349729213332
$ws.Range("I19").Value = 0.7402018593529966
$ws.Range("L19").Value = 0.2681539931701877
$ws.Range("M19").Value = 0.252694082475216
$ws.Range("N19").Value = 1.606547036629536
$ws.Range("O19").Value = 2.055126252844843
# Row 20
$ws.Range("B20").Value = 0.9324146719220892
$ws.Range("D20").Value = 0.00546662080218141
$ws.Range("E20").Value = 0.526536304640679
$ws.Range("F20").Value = 0.6601218633263102
$ws.Range("G20").Value = 0.002357196563635427
$ws.Range("I20").Value = 0.7359787154145967
$ws.Range("L20").Value = 0.2756316824651464
$ws.Range("M20").Value = 0.2586535447795484
$ws.Range("N20").Value = 1.611337307500762
$ws.Range("O20").Value = 2.085548608991587
# Row 21
$ws.Range("B21").Value = 0.996214819414206
$ws.Range("D21").Value = 0.005953982085944887
$ws.Range("E21").Value = 0.5448569202122187
$ws.Range("F21").Value = 0.6963015311832095
$ws.Range("G21").Value = 0.002353449605153928
$ws.Range("I21").Value = 0.7223031082615465
$ws.Range("L21").Value = 0.3008056320326773
$ws.Range("M21").Value = 0.2787237642967924
$ws.Range("N21").Value = 1.628077577633945
$ws.Range("O21").Value = 2.189405936209141
# Row 22
$ws.Range("B22").Value = 1.038028072324664
$ws.Range("D22").Value = 0.006271338881873589
$ws.Range("E22").Value = 0.5568080075980788
$ws.Range("F22").Value = 0.7203203669404843
$ws.Range("G22").Value = 0.002351095386738709
$ws.Range("I22").Value = 0.7137483170892214
$ws.Range("L22").Value = 0.3172942042718603
$ws.Range("M22").Value = 0.2918745975723169
$ws.Range("N22").Value = 1.639495599053134
$ws.Range("O22").Value = 2.258511877538524
# Row 23
$ws.Range("B23").Value = 1.015701263357414
$ws.Range("D23").Value = 0.006102067924562249
$ws.Range("E23").Value = 0.5504316620433016
$ws.Range("F23").Value = 0.7074668131167243
$ws.Range("G23").Value = 0.002352343244505443
$ws.Range("I23").Value = 0.7182791098791554
$ws.Range("L23").Value = 0.308490736389615
$ws.Range("M23").Value = 0.2848527505419085
$ws.Range("N23").Value = 1.633358163534496
$ws.Range("O23").Value = 2.221515940219433
# Row 24
$ws.Range("B24").Value = 0.9313739287644864
$ws.Range("D24").Value = 0.005458635969272763
$ws.Range("E24").Value = 0.5262365099676884
$ws.Range("F24").Value = 0.6595368080229207
$ws.Range("G24").Value = 0.002357259409115181
$ws.Range("I24").Value = 0.736208698792693
$ws.Range("L24").Value = 0.2752208523683493
$ws.Range("M24").Value = 0.2583260982179851
$ws.Range("N24").Value = 1.61107180989579
$ws.Range("O24").Value = 2.083871792828973
# Row 25
$ws.Range("B25").Value = 0.8409460971047338
$ws.Range("D25").Value = 0.004759803926578599
$ws.Range("E25").Value = 0.5000532517700336
$ws.Range("F25").Value = 0.6094301192209457
$ws.Range("G25").Value = 0.002362971257204787
$ws.Range("I25").Value = 0.7571916079001351
$ws.Range("L25").Value = 0.239496358570932
$ws.Range("M25").Value = 0.2298669994567533
$ws.Range("N25").Value = 1.589106791721791
$ws.Range("O25").Value = 1.940643935681805

Write-Host "Updated 240 cells"